$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet1: record a "Widthdrawn" count on row 3, and bump the LastDate
# timestamp to match the later run.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("D3").Value = 2
$ws1.Range("E3").Value = 44830.43886574074

# ---------------------------------------------------------------------
# Sheet2: new "withdraw" log - add DataTimestamp / Widthdrawn columns
# and append the rows of people who were messaged / connected with.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("F1").Value = "DataTimestamp"
$ws2.Range("G1").Value = "Widthdrawn"

$ws2.Range("A2").Value = "Nikhita Singh"
$ws2.Range("B2").Value = "talent acquisition"
$ws2.Range("C2").Value = "Bengaluru"
$ws2.Range("D2").Value = "Message - Sent"
$ws2.Range("E2").Value = "https://www.linkedin.com/in/nikhita-singh-a8423878?miniProfileUrn=urn%3Ali%3Afs_miniProfile%3AACoAABBohwQBVxAaxnlHGcEhGYPLbUk2K4cfspM"
$ws2.Range("F2").Value = 44830.40488425926

$ws2.Range("A3").Value = "preethi kumar"
$ws2.Range("B3").Value = "talent acquisition"
$ws2.Range("C3").Value = "Bengaluru"
$ws2.Range("D3").Value = "Message - Already Sent"
$ws2.Range("E3").Value = "https://www.linkedin.com/in/preethi-kumar-533603145?miniProfileUrn=urn%3Ali%3Afs_miniProfile%3AACoAACMmrzcBuZT5gSFrQPxrZBxQK3goTwkZbDc"
$ws2.Range("F3").Value = 44830.404988425929

$ws2.Range("A4").Value = "Aishwarya BP"
$ws2.Range("B4").Value = "talent acquisition"
$ws2.Range("C4").Value = "Bengaluru"
$ws2.Range("D4").Value = "Message - Sent"
$ws2.Range("E4").Value = "https://www.linkedin.com/in/aishwarya-bp-37a059103?miniProfileUrn=urn%3Ali%3Afs_miniProfile%3AACoAABorg84BMXvgvbAyq0JxoxyIAgk8RahvhDQ"
$ws2.Range("F4").Value = 44830.40525462963

$ws2.Range("A5").Value = "Shruthi Abbar"
$ws2.Range("B5").Value = "talent acquisition"
$ws2.Range("C5").Value = "Bengaluru"
$ws2.Range("D5").Value = "Connect - Sent"
$ws2.Range("E5").Value = "https://www.linkedin.com/in/shruthi-abbar-47685a177?miniProfileUrn=urn%3Ali%3Afs_miniProfile%3AACoAACnzaxoBySAvXwpC8YnbJwBuTWv96NXm6To"
$ws2.Range("F5").Value = 44830.405300925922

$ws2.Range("A6").Value = "RIVIN MATHEW"
$ws2.Range("B6").Value = "talent acquisition"
$ws2.Range("C6").Value = "Bengaluru"
$ws2.Range("D6").Value = "Connect - Sent"
$ws2.Range("E6").Value = "https://www.linkedin.com/in/rivin-mathew-0a64b0135?miniProfileUrn=urn%3Ali%3Afs_miniProfile%3AACoAACD_u34B2Nqz595efGkVi1Nj_9ZsZo8c0TI"
$ws2.Range("F6").Value = 44830.405347222222

$ws2.Range("A7").Value = "Aarthi Raju"
$ws2.Range("B7").Value = "talent acquisition"
$ws2.Range("C7").Value = "Bengaluru"
$ws2.Range("D7").Value = "Message - Sent"
$ws2.Range("E7").Value = "https://www.linkedin.com/in/aarthi-raju-792a0b197?miniProfileUrn=urn%3Ali%3Afs_miniProfile%3AACoAAC4-9p4BTu6bJ_ft8aDkNyJyZvqS2H55-M4"
$ws2.Range("F7").Value = 44830.437939814816

$ws2.Range("A8").Value = "Smitha H"
$ws2.Range("B8").Value = "talent acquisition"
$ws2.Range("C8").Value = "Bengaluru"
$ws2.Range("D8").Value = "Message - Sent"
$ws2.Range("E8").Value = "https://www.linkedin.com/in/smitha-h-3672a2155?miniProfileUrn=urn%3Ali%3Afs_miniProfile%3AACoAACU9YX8BPJ8Jtjzq1onul_LeeohCQskIWTU"
$ws2.Range("F8").Value = 44830.438217592593

$ws2.Range("A9").Value = "Sonal Ranjit"
$ws2.Range("B9").Value = "talent acquisition"
$ws2.Range("C9").Value = "Bengaluru"
$ws2.Range("D9").Value = "Message - Sent"
$ws2.Range("E9").Value = "https://www.linkedin.com/in/sonal-ranjit-a77ba8164?miniProfileUrn=urn%3Ali%3Afs_miniProfile%3AACoAACdU6fsBWhp9QE47Xy_kT6iQQwHvNpr27b4"
$ws2.Range("F9").Value = 44830.43849537037

$ws2.Range("A10").Value = "Parul Narayan"
$ws2.Range("B10").Value = "talent acquisition"
$ws2.Range("C10").Value = "Bengaluru"
$ws2.Range("D10").Value = "Message - Sent"
$ws2.Range("E10").Value = "https://www.linkedin.com/in/parul-narayan-b551391a4?miniProfileUrn=urn%3Ali%3Afs_miniProfile%3AACoAAC-3wjEB_vXizWzTrx2OKoZwG2m9aZoFMNw"
$ws2.Range("F10").Value = 44830.438761574071

$ws2.Range("A11").Value = "Natasha Castelino"
$ws2.Range("B11").Value = "talent acquisition"
$ws2.Range("C11").Value = "Bengaluru"
$ws2.Range("D11").Value = "Connect - Sent"
$ws2.Range("E11").Value = "https://www.linkedin.com/in/natasha-castelino-92590b7b?miniProfileUrn=urn%3Ali%3Afs_miniProfile%3AACoAABEK620B9eP5H9iv543swL-IiLJKb-ALHMo"
$ws2.Range("F11").Value = 44830.438807870371

# Re-use the existing date/time number format (the same style already
# applied to Sheet1's LastDate column) for the new timestamp column,
# instead of letting a brand-new style get created.
$ws1.Range("E3").Copy()
$ws2.Range("F2:F11").PasteSpecial(-4122)

# Size the two new columns to fit their content, like the rest of the
# sheet's columns.
$ws2.Columns.Item(6).AutoFit()
$ws2.Columns.Item(7).AutoFit()

# ---------------------------------------------------------------------
# Selection / active-sheet bookkeeping: the author ended up back on
# Sheet1 with B4 selected, leaving Sheet2's cursor parked at G2.
# ---------------------------------------------------------------------
$ws2.Range("G2").Select()
$ws1.Activate()
$ws1.Range("B4").Select()
